# Apply updated NATMI metrics (ligand/receptor-expressing cell counts 1 -> 3 per Dr Hou advice)
# and the resulting recomputed derived values, per worksheet cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E2" = 3
    "G2" = 2.098888
    "H2" = 6.296664
    "I2" = 0.1082453658858517
    "J2" = 0.1082453658858517
    "K2" = 3
    "M2" = 23.59622066666667
    "N2" = 70.788662
    "O2" = 0.6996728317814862
    "P2" = 0.6996728317814862
    "Q2" = 49.52582440261867
    "R2" = 445.732419623568
    "S2" = 0.07573634167657693
    "T2" = 0.0757363416765769
    "E3" = 3
    "G3" = 2.098888
    "H3" = 6.296664
    "I3" = 0.1082453658858517
    "J3" = 0.1082453658858517
    "K3" = 3
    "M3" = 7.778025666666667
    "N3" = 23.334077
    "O3" = 0.2306332577891816
    "P3" = 0.2306332577891816
    "Q3" = 16.32520473545867
    "R3" = 146.926842619128
    "S3" = 0.02496498137483592
    "T3" = 0.02496498137483591
    "E4" = 3
    "G4" = 2.098888
    "H4" = 6.296664
    "I4" = 0.1082453658858517
    "J4" = 0.1082453658858517
    "K4" = 3
    "M4" = 2.350402666666667
    "N4" = 7.051208000000001
    "O4" = 0.06969391042933218
    "P4" = 0.06969391042933218
    "Q4" = 4.933231952234668
    "R4" = 44.399087570112
    "S4" = 0.007544042834438836
    "T4" = 0.007544042834438834
    "E5" = 3
    "G5" = 15.87514366666667
    "H5" = 47.625431
    "I5" = 0.8187243600843848
    "J5" = 0.8187243600843847
    "K5" = 3
    "M5" = 23.59622066666667
    "N5" = 70.788662
    "O5" = 0.6996728317814862
    "P5" = 0.6996728317814862
    "Q5" = 374.5933930737024
    "R5" = 3371.340537663322
    "S5" = 0.5728391914687266
    "T5" = 0.5728391914687266
    "E6" = 3
    "G6" = 15.87514366666667
    "H6" = 47.625431
    "I6" = 0.8187243600843848
    "J6" = 0.8187243600843847
    "K6" = 3
    "M6" = 7.778025666666667
    "N6" = 23.334077
    "O6" = 0.2306332577891816
    "P6" = 0.2306332577891816
    "Q6" = 123.4772749013541
    "R6" = 1111.295474112187
    "S6" = 0.1888250663976247
    "T6" = 0.1888250663976247
    "E7" = 3
    "G7" = 15.87514366666667
    "H7" = 47.625431
    "I7" = 0.8187243600843848
    "J7" = 0.8187243600843847
    "K7" = 3
    "M7" = 2.350402666666667
    "N7" = 7.051208000000001
    "O7" = 0.06969391042933218
    "P7" = 0.06969391042933218
    "Q7" = 37.31298000784978
    "R7" = 335.816820070648
    "S7" = 0.05706010221803342
    "T7" = 0.05706010221803341
    "E8" = 3
    "G8" = 1.416064
    "H8" = 4.248192
    "I8" = 0.07303027402976368
    "J8" = 0.07303027402976367
    "K8" = 3
    "M8" = 23.59622066666667
    "N8" = 70.788662
    "O8" = 0.6996728317814862
    "P8" = 0.6996728317814862
    "Q8" = 33.41375862212266
    "R8" = 300.723827599104
    "S8" = 0.05109729863618269
    "T8" = 0.05109729863618268
    "E9" = 3
    "G9" = 1.416064
    "H9" = 4.248192
    "I9" = 0.07303027402976368
    "J9" = 0.07303027402976367
    "K9" = 3
    "M9" = 7.778025666666667
    "N9" = 23.334077
    "O9" = 0.2306332577891816
    "P9" = 0.2306332577891816
    "Q9" = 11.01418213764267
    "R9" = 99.12763923878398
    "S9" = 0.01684321001672107
    "T9" = 0.01684321001672106
    "E10" = 3
    "G10" = 1.416064
    "H10" = 4.248192
    "I10" = 0.07303027402976368
    "J10" = 0.07303027402976367
    "K10" = 3
    "M10" = 2.350402666666667
    "N10" = 7.051208000000001
    "O10" = 0.06969391042933218
    "P10" = 0.06969391042933218
    "Q10" = 3.328320601770667
    "R10" = 29.954885415936
    "S10" = 0.005089765376859934
    "T10" = 0.005089765376859933
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

